$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column E (Description instead of Short Description)
$ws.Range("E1").Value = "Description"

# Update speaker(s) column F with markdown-style links
$ws.Range("F2").Value = "[Lauren Chenarides](https://wpcarey.asu.edu/people/profile/3153292), [Drew Hanks](https://ehe.osu.edu/human-sciences/directory?id=hanks.46)"
$ws.Range("F3").Value = "[Tim Beatty](https://are.ucdavis.edu/people/faculty/tim-beatty/)"
$ws.Range("F4").Value = "[Julia Lane](https://wagner.nyu.edu/community/faculty/julia-lane)"
$ws.Range("F5").Value = "[Julia Lane](https://wagner.nyu.edu/community/faculty/julia-lane), [Jason Owen-Smith](https://lsa.umich.edu/soc/people/faculty/jdos.html)"
$ws.Range("F6").Value = "[Jason Owen-Smith](https://lsa.umich.edu/soc/people/faculty/jdos.html)"
$ws.Range("F8").Value = "[Andi Carlson](https://www.ers.usda.gov/authors/ers-staff-directory/andrea-carlson/)"
$ws.Range("F9").Value = "[Ayaz Hyder](https://cph.osu.edu/people/ahyder), [Charlotte Ambrozek](https://are.ucdavis.edu/people/grad-students/phd/charlotte-ambrozek/)"
$ws.Range("F10").Value = "[Lauren Chenarides](https://wpcarey.asu.edu/people/profile/3153292), [Drew Hanks](https://ehe.osu.edu/human-sciences/directory?id=hanks.46)"

# Update the active selection to match the saved view state
$ws.Range("E18").Select()
